$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "展览"  (Exhibition)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Insert a brand-new row at position 43 (pushes the old row 43 down to 44).
$ws1.Rows.Item(43).Insert()

# The insert leaves the new row (and sometimes the shifted row) with a
# slightly different auto-generated style than the rest of the table, so
# re-stamp both rows' formatting from the known-good row directly above.
$ws1.Range("A42:I42").Copy()
$ws1.Range("A43:I44").PasteSpecial(-4122)

# Column A is a simple sequential row index (0,1,2,...); restore it since the
# plain row-insert leaves A43 blank and A44 holding the old (pre-shift) value.
$ws1.Range("A43").Value = 42
$ws1.Range("A44").Value = 43

# Populate the new event row 43.
$ws1.Range("B43").NumberFormat = "@"
$ws1.Range("B43").Value = "2024-09-17"
$ws1.Range("B43").Style = "Normal"
$ws1.Range("C43").Value = "广州·樱漫潮玩动漫嘉年华"
$ws1.Range("D43").Value = "开创大道2666号 宝能国际体育演艺中心"
$ws1.Range("E43").Value = "2024.09.17 10:00-09.17 17:00"
$ws1.Range("F43").Value = 0
$ws1.Range("G43").Value = 9.9
$ws1.Range("H43").Value = "https://show.bilibili.com/platform/detail.html?id=88936"
$ws1.Range("I43").Value = "//i0.hdslb.com/bfs/openplatform/202407/v2XniRgx1719799991150.jpeg"

# "想去人数" (want-to-go count) refresh across the rest of the sheet.
$ws1.Range("F2").Value = 277
$ws1.Range("F3").Value = 0
$ws1.Range("F5").Value = 866
$ws1.Range("F6").Value = 476
$ws1.Range("F7").Value = 0
$ws1.Range("G9").Value = "暂时售罄"
$ws1.Range("F12").Value = 515
$ws1.Range("F13").Value = 672
$ws1.Range("F17").Value = 180
$ws1.Range("F19").Value = 37
$ws1.Range("F20").Value = 79
$ws1.Range("F24").Value = 0
$ws1.Range("F25").Value = 524
$ws1.Range("F26").Value = 371
$ws1.Range("F27").Value = 522
$ws1.Range("F28").Value = 0
$ws1.Range("F29").Value = 31
$ws1.Range("F30").Value = 0
$ws1.Range("F32").Value = 0
$ws1.Range("F35").Value = 169
$ws1.Range("F36").Value = 0
$ws1.Range("F37").Value = 0
$ws1.Range("F40").Value = 0
$ws1.Range("F41").Value = 319
$ws1.Range("F42").Value = 70

# ---------------------------------------------------------------------------
# Sheet 2: "演出"  (Performance)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

# Insert a brand-new row at position 19 (pushes the old row 19 down to 20).
$ws2.Rows.Item(19).Insert()
$ws2.Range("A18:I18").Copy()
$ws2.Range("A19:I20").PasteSpecial(-4122)

# Restore the sequential row index in column A.
$ws2.Range("A19").Value = 18
$ws2.Range("A20").Value = 19

# Populate the new event row 19.
$ws2.Range("B19").NumberFormat = "@"
$ws2.Range("B19").Value = "2024-09-22"
$ws2.Range("B19").Style = "Normal"
$ws2.Range("C19").Value = "广州·VGL中国巡演 2024 VIDEO GAMES LIVE 魔兽世界音乐会"
$ws2.Range("D19").Value = "东风中路299号 广州中山纪念堂"
$ws2.Range("E19").Value = "2024.09.22 19:30-09.22 21:10"
$ws2.Range("F19").Value = 1
$ws2.Range("G19").Value = 180
$ws2.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=88919"
$ws2.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202407/Bn0OQ6ef1720428966541.jpeg"

# The shifted-down row (now row 20) also picks up a refreshed want-to-go count.
$ws2.Range("F20").Value = 5

# "想去人数" refresh across the rest of the sheet.
$ws2.Range("F4").Value = 0
$ws2.Range("F5").Value = 4354
$ws2.Range("F7").Value = 0
$ws2.Range("F9").Value = 5
$ws2.Range("F10").Value = 69
$ws2.Range("F12").Value = 0
$ws2.Range("F14").Value = 67
$ws2.Range("F15").Value = 158
$ws2.Range("F16").Value = 0
$ws2.Range("F17").Value = 4351
$ws2.Range("F18").Value = 0

# ---------------------------------------------------------------------------
# Sheet 3: "本地生活"  (Local life)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Range("F2").Value = 0
$ws3.Range("F3").Value = 418
$ws3.Range("F4").Value = 238

# ---------------------------------------------------------------------------
# Sheet 4: "全部类型"  (All types)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value = 418
$ws4.Range("F4").Value = 238
$ws4.Range("F6").Value = 928
$ws4.Range("F7").Value = 87
$ws4.Range("F8").Value = 866
$ws4.Range("F10").Value = 0
$ws4.Range("F11").Value = 38690
$ws4.Range("G11").Value = "暂时售罄"
$ws4.Range("F13").Value = 321
$ws4.Range("F16").Value = 7
$ws4.Range("F17").Value = 0
$ws4.Range("F18").Value = 0
$ws4.Range("F19").Value = 672
$ws4.Range("F20").Value = 0
$ws4.Range("F21").Value = 81
$ws4.Range("F22").Value = 81
$ws4.Range("F23").Value = 0
$ws4.Range("F24").Value = 180
$ws4.Range("F26").Value = 6
$ws4.Range("F29").Value = 985
$ws4.Range("F31").Value = 524
$ws4.Range("F32").Value = 371
$ws4.Range("F36").Value = 4
$ws4.Range("F39").Value = 804
$ws4.Range("F40").Value = 0
$ws4.Range("F42").Value = 208
$ws4.Range("F43").Value = 2

Write-Host "edit complete"
